$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the affected rows.
# These values represent a "repull" of data / mean calculation that changed
# the stored dSF figures while leaving dS0 (column E) and everything else intact.
$updates = @{
    2  = -2
    3  = -1
    4  = 1
    5  = 3
    6  = 4
    9  = -4
    10 = -4
    11 = 3
    12 = 3
    13 = -2
    14 = 3
    16 = -2
    17 = 3
    18 = 4
    19 = -2
    20 = 2
    22 = -1
    24 = 3
    25 = -1
    26 = -1
    27 = 2
    29 = 2
    30 = 3
    31 = -1
    32 = -2
    33 = 4
    35 = -5
    36 = 2
    37 = -1
    38 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
